$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place substring replace to preserve surrounding runs) ---
$ws.Range("A8").Characters(21, 2).Text = "38"
$ws.Range("C9").Characters(46, 9).Text = "9/22/2024"
$ws.Range("C9").Characters(27, 8).Text = "9/16/2024"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -33.333333333333
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 9
$ws.Range("H15").Value = 122.222222222222
$ws.Range("I15").Value = 117
$ws.Range("J15").Value = 85
$ws.Range("K15").Value = 37.647058823529
$ws.Range("L15").Value = -11.363636363636
$ws.Range("M15").Value = 62.5
$ws.Range("N15").Value = -30.357142857142
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 41
$ws.Range("E16").Value = -19.512195121951
$ws.Range("F16").Value = 143
$ws.Range("G16").Value = 139
$ws.Range("H16").Value = 2.877697841726
$ws.Range("I16").Value = 1219
$ws.Range("J16").Value = 1325
$ws.Range("K16").Value = -8
$ws.Range("L16").Value = -20.534550195567
$ws.Range("M16").Value = 33.662280701754
$ws.Range("N16").Value = -84.187313529640
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = 27.777777777777
$ws.Range("F17").Value = 184
$ws.Range("H17").Value = 12.883435582822
$ws.Range("I17").Value = 1715
$ws.Range("J17").Value = 1611
$ws.Range("K17").Value = 6.455617628801
$ws.Range("L17").Value = 10.859728506787
$ws.Range("M17").Value = 88.254665203073
$ws.Range("N17").Value = -32.665881429132
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = -17.073170731707
$ws.Range("F18").Value = 122
$ws.Range("G18").Value = 166
$ws.Range("H18").Value = -26.506024096385
$ws.Range("I18").Value = 1335
$ws.Range("J18").Value = 1543
$ws.Range("K18").Value = -13.480233311730
$ws.Range("L18").Value = -39.290586630286
$ws.Range("M18").Value = 2.771362586605
$ws.Range("N18").Value = -84.672789896670
$ws.Range("C19").Value = 195
$ws.Range("D19").Value = 243
$ws.Range("E19").Value = -19.753086419753
$ws.Range("F19").Value = 833
$ws.Range("G19").Value = 930
$ws.Range("H19").Value = -10.430107526881
$ws.Range("I19").Value = 7534
$ws.Range("J19").Value = 8409
$ws.Range("K19").Value = -10.405517897490
$ws.Range("L19").Value = -9.664268585131
$ws.Range("M19").Value = 0.119601328903
$ws.Range("N19").Value = -68.387042631755
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 16
$ws.Range("E20").Value = -43.75
$ws.Range("F20").Value = 40
$ws.Range("G20").Value = 70
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 333
$ws.Range("J20").Value = 470
$ws.Range("K20").Value = -29.148936170212
$ws.Range("L20").Value = -34.319526627218
$ws.Range("M20").Value = 14.432989690721
$ws.Range("N20").Value = -92.515171948752
$ws.Range("C21").Value = 321
$ws.Range("D21").Value = 383
$ws.Range("E21").Value = -16.187989556135
$ws.Range("F21").Value = 1343
$ws.Range("G21").Value = 1477
$ws.Range("H21").Value = -9.072444143534
$ws.Range("I21").Value = 12269
$ws.Range("J21").Value = 13457
$ws.Range("K21").Value = -8.828119194471
$ws.Range("L21").Value = -14.112705635281
$ws.Range("M21").Value = 11.333938294010
$ws.Range("N21").Value = -74.158557647752
$ws.Range("C22").Value = 10
$ws.Range("E22").Value = -23.076923076923
$ws.Range("F22").Value = 36
$ws.Range("G22").Value = 48
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 439
$ws.Range("J22").Value = 477
$ws.Range("K22").Value = -7.966457023060
$ws.Range("L22").Value = -8.541666666666
$ws.Range("M22").Value = 14.025974025974
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = 17.857142857142
$ws.Range("I23").Value = 289
$ws.Range("J23").Value = 293
$ws.Range("K23").Value = -1.365187713310
$ws.Range("L23").Value = -13.731343283582
$ws.Range("M23").Value = 11.583011583011
$ws.Range("C24").Value = 493
$ws.Range("D24").Value = 369
$ws.Range("E24").Value = 33.604336043360
$ws.Range("F24").Value = 1814
$ws.Range("G24").Value = 1643
$ws.Range("H24").Value = 10.407790626902
$ws.Range("I24").Value = 16172
$ws.Range("J24").Value = 15060
$ws.Range("K24").Value = 7.383798140770
$ws.Range("L24").Value = 1.512773837172
$ws.Range("M24").Value = 31.919406150583
$ws.Range("C25").Value = 392
$ws.Range("D25").Value = 312
$ws.Range("E25").Value = 25.641025641025
$ws.Range("F25").Value = 1455
$ws.Range("G25").Value = 1298
$ws.Range("H25").Value = 12.095531587057
$ws.Range("I25").Value = 13639
$ws.Range("J25").Value = 12474
$ws.Range("K25").Value = 9.339426006092
$ws.Range("L25").Value = 0.701417601890
$ws.Range("C26").Value = 130
$ws.Range("D26").Value = 82
$ws.Range("E26").Value = 58.536585365853
$ws.Range("F26").Value = 441
$ws.Range("G26").Value = 393
$ws.Range("H26").Value = 12.213740458015
$ws.Range("I26").Value = 3688
$ws.Range("J26").Value = 3647
$ws.Range("K26").Value = 1.124211680833
$ws.Range("L26").Value = 8.279506752789
$ws.Range("M26").Value = 40.334855403348
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 15
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 177
$ws.Range("J27").Value = 152
$ws.Range("K27").Value = 16.447368421052
$ws.Range("L27").Value = -14.077669902912
$ws.Range("C28").Value = 25
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 56.25
$ws.Range("F28").Value = 101
$ws.Range("G28").Value = 61
$ws.Range("H28").Value = 65.573770491803
$ws.Range("I28").Value = 711
$ws.Range("J28").Value = 654
$ws.Range("K28").Value = 8.715596330275
$ws.Range("L28").Value = -0.280504908835
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 6
$ws.Range("I29").Value = 30
$ws.Range("K29").Value = 3.448275862068
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = 15.384615384615
$ws.Range("N29").Value = -70
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 6
$ws.Range("I30").Value = 26
$ws.Range("K30").Value = 8.333333333333
$ws.Range("L30").Value = -25.714285714285
$ws.Range("M30").Value = 30
$ws.Range("N30").Value = -70.454545454545
$ws.Range("G31").Value = 10
$ws.Range("L31").Value = -20.161290322580
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 12
$ws.Range("K33").Value = 33.333333333333

# --- Text-to-text updates (non-numeric-looking replacement strings; style unchanged) ---
$ws.Range("D14").Value = "***.*"

# --- Numeric cells becoming text placeholders ("0" / "***.*"); target style = A14's style ---
$ws.Range("C14").Formula = "=""0"""
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("G29").Formula = "=""0"""
$ws.Range("G29").Copy()
$ws.Range("G29").PasteSpecial(-4163)
$ws.Range("H29").Formula = "=""***.*"""
$ws.Range("H29").Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("G30").Formula = "=""0"""
$ws.Range("G30").Copy()
$ws.Range("G30").PasteSpecial(-4163)
$ws.Range("H30").Formula = "=""***.*"""
$ws.Range("H30").Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=""0"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = "=""***.*"""
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("F33").Formula = "=""0"""
$ws.Range("F33").Copy()
$ws.Range("F33").PasteSpecial(-4163)
$ws.Range("A14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("F33").PasteSpecial(-4122)

# --- Text placeholder cells becoming numeric values; restore numeric style from same-row peers ---
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("H33").Copy()
$ws.Range("E33").PasteSpecial(-4122)

$excel.CutCopyMode = $false
